$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1.xml)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 24
$ws.Range("B3").Value = 24
$ws.Range("E3").Value = 39
$ws.Range("I3").Value = 46
$ws.Range("C6").Value = 118
$ws.Range("D6").Value = 120
$ws.Range("F6").Value = 150
$ws.Range("G6").Value = 159
$ws.Range("I6").Value = 145
$ws.Range("J6").Value = 110
$ws.Range("B7").Value = 152
$ws.Range("C7").Value = 161
$ws.Range("D7").Value = 181
$ws.Range("E7").Value = 184
$ws.Range("F7").Value = 210
$ws.Range("G7").Value = 219
$ws.Range("I7").Value = 219
$ws.Range("J7").Value = 209

# By Neighborhood (sheet2.xml)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 4
$ws.Range("C7").Value = 15
$ws.Range("G7").Value = 15
$ws.Range("J24").Value = 6
$ws.Range("B33").Value = 10
$ws.Range("E33").Value = 8
$ws.Range("I33").Value = 14
$ws.Range("D50").Value = 28
$ws.Range("E50").Value = 30
$ws.Range("G50").Value = 29
$ws.Range("I50").Value = 36
$ws.Range("F60").Value = 2
$ws.Range("G72").Value = 5
$ws.Range("I73").Value = 9
$ws.Range("J73").Value = 4
$ws.Range("C74").Value = 4
$ws.Range("F74").Value = 5
$ws.Range("B92").Value = 152
$ws.Range("C92").Value = 161
$ws.Range("D92").Value = 181
$ws.Range("E92").Value = 184
$ws.Range("F92").Value = 210
$ws.Range("G92").Value = 219
$ws.Range("I92").Value = 219
$ws.Range("J92").Value = 209

# Roseland (sheet4.xml)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 5

# Austin (sheet7.xml)
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("C5").Value = 12
$ws.Range("G5").Value = 13
$ws.Range("C6").Value = 15
$ws.Range("G6").Value = 15

# Grand Crossing (sheet11.xml)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("B3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("I3").Value = 6
$ws.Range("B6").Value = 10
$ws.Range("E6").Value = 8
$ws.Range("I6").Value = 14

# Armour Square (sheet12.xml)
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 4

# Rush & Division (sheet20.xml)
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("H3").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("H5").Value = 9
$ws.Range("I5").Value = 4

# Loop (sheet22.xml)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E2").Value = 2
$ws.Range("D6").Value = 16
$ws.Range("G6").Value = 19
$ws.Range("I6").Value = 26
$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 30
$ws.Range("G7").Value = 29
$ws.Range("I7").Value = 36

# Sheffield & DePaul (sheet28.xml)
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("F2").Value = 1
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 4
$ws.Range("F6").Value = 5

# New City (sheet43.xml)
$ws = $wb.Worksheets.Item('New City')
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 2

# Edgewater (sheet44.xml)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 6
